$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# 1) Rename the header row: the "*_old" columns become "*_FV2310" and the
#    "*_new" columns become "*_FV2404" (column K stays "diff").
# -----------------------------------------------------------------------
$ws.Range("A1").Value = "Segmentname_FV2310"
$ws.Range("B1").Value = "Segmentgruppe_FV2310"
$ws.Range("C1").Value = "Segment_FV2310"
$ws.Range("D1").Value = "Datenelement_FV2310"
$ws.Range("E1").Value = "Segment ID_FV2310"
$ws.Range("F1").Value = "Code_FV2310"
$ws.Range("G1").Value = "Qualifier_FV2310"
$ws.Range("H1").Value = "Beschreibung_FV2310"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("J1").Value = "Bedingung_FV2310"

$ws.Range("L1").Value = "Segmentname_FV2404"
$ws.Range("M1").Value = "Segmentgruppe_FV2404"
$ws.Range("N1").Value = "Segment_FV2404"
$ws.Range("O1").Value = "Datenelement_FV2404"
$ws.Range("P1").Value = "Segment ID_FV2404"
$ws.Range("Q1").Value = "Code_FV2404"
$ws.Range("R1").Value = "Qualifier_FV2404"
$ws.Range("S1").Value = "Beschreibung_FV2404"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("U1").Value = "Bedingung_FV2404"

# -----------------------------------------------------------------------
# 2) Freeze the header row (split under row 1, top-left cell of the
#    scrolling pane is A2).
# -----------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# -----------------------------------------------------------------------
# 3) Turn A1:U68 into a native Excel Table ("Table1") with an AutoFilter.
#    The header row already carries bold/fill/border/center/wrap
#    formatting (style index 1); temporarily reset it to the default
#    "Normal" style before creating the table so that Excel does not
#    capture that look as a one-off table header style (dxf), then put
#    the same formatting back afterwards so the visible result - and the
#    underlying style index used by every header cell - is unchanged.
# -----------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$headerRange.Style = "Normal"

$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U68"), 0, 1)
$tbl.Name = "Table1"

$headerRange.Font.Bold = $true
$headerRange.Interior.Color = 14277081
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.WrapText = $true
